# Actualización automática 2025-08-27 10:00:10
#
# Updates August ("agosto") sales figures for several clients of
# GUERRERO FAREZ FABIAN MAURICIO across the three report sheets:
#   - "VENTAS POR GRUPO"      : sales broken down by product group
#   - "VENTA MENSUAL"         : sales broken down by month (column F = agosto)
#   - "CUMPLIMIENTO MENSUAL"  : budget vs. sales compliance summary by group

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# AGUIMPORT-AGUILAR IMPORTACIONES S.A.S. (row 5) - PORCELANATO
$wsGrupo.Range("M5").Value = 2587.12

# BUSTAMANTE ROSERO MARCO TULIO (row 12) - PIEDRA SINTERIZADA
$wsGrupo.Range("L12").Value = 760.3200000000001

# FEIJOO MARIN MAURICIO ENRIQUE (row 22) - PORCELANATO
$wsGrupo.Range("M22").Value = 620.21

# SALAS NOBLECILLA MARIA SUSANA (row 48) - INODOROS, PORCELANATO, PUERTAS DE SEGURIDAD
$wsGrupo.Range("H48").Value = 447.3
$wsGrupo.Range("M48").Value = 366.83
$wsGrupo.Range("N48").Value = 248.48

# WONG SANCHEZ CLAUDIA PAULINA (row 52) - PORCELANATO
$wsGrupo.Range("M52").Value = 1185.28

# WONG SANCHEZ PAULA SOFIA (row 53) - SAL SOLUBLE
$wsGrupo.Range("O53").Value = 4253.64

# Row 55 totals ("X de 53" clients with sales in that group)
$wsGrupo.Range("L55").Value = "6 de 53"
$wsGrupo.Range("M55").Value = "17 de 53"
$wsGrupo.Range("N55").Value = "1 de 53"
$wsGrupo.Range("O55").Value = "7 de 53"

# ---------------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL (column F = agosto)
# ---------------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F5").Value = 4900.9
$wsMensual.Range("F12").Value = 2601.27
$wsMensual.Range("F22").Value = 6477.99
$wsMensual.Range("F48").Value = 1204.72
$wsMensual.Range("F52").Value = 2471.15
$wsMensual.Range("F53").Value = 4253.64
$wsMensual.Range("F55").Value = 79717.36

# ---------------------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL (PRESUPUESTO / VENTA / POR CUMPLIR / CUMPLIMIENTO)
# ---------------------------------------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# INODOROS (row 7)
$wsCumplimiento.Range("D7").Value = 1788.3
$wsCumplimiento.Range("E7").Value = 611.7
$wsCumplimiento.Range("F7").Value = 0.7451249999999999

# PIEDRA SINTERIZADA (row 15)
$wsCumplimiento.Range("D15").Value = 6923.61
$wsCumplimiento.Range("E15").Value = 6576.39
$wsCumplimiento.Range("F15").Value = 0.51286

# PORCELANATO (row 16)
$wsCumplimiento.Range("D16").Value = 34672.62
$wsCumplimiento.Range("E16").Value = 21387.07999999999
$wsCumplimiento.Range("F16").Value = 0.6184945691824966

# PUERTAS DE SEGURIDAD (row 17)
$wsCumplimiento.Range("D17").Value = 248.48
$wsCumplimiento.Range("E17").Value = 435.52
$wsCumplimiento.Range("F17").Value = 0.3632748538011696

# SAL SOLUBLE (row 18)
$wsCumplimiento.Range("D18").Value = 11382.06
$wsCumplimiento.Range("E18").Value = -8182.059999999999
$wsCumplimiento.Range("F18").Value = 3.55689375

# TOTAL (row 19)
$wsCumplimiento.Range("D19").Value = 79717.36
$wsCumplimiento.Range("E19").Value = 37722.33064517914
$wsCumplimiento.Range("F19").Value = 0.6787940223791143
